# Update the "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated stats, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1392
    $ws.Range("F3").Value = 2659
    $ws.Range("F4").Value = 546
    $ws.Range("F5").Value = 88
    $ws.Range("F6").Value = 6612
    $ws.Range("F7").Value = 471
    $ws.Range("F8").Value = 9
    $ws.Range("F9").Value = 10
    $ws.Range("F10").Value = 51
    $ws.Range("F12").Value = 132
}
